# edit.ps1 -- apply documented change to The_BOIS.docx
#
# The diff does two things to the last paragraph of the body ("A movie
# themed application ..."):
#   1. Removes the <w:proofErr w:type="gramStart"/> / <w:proofErr
#      w:type="gramEnd"/> grammar-check markers that wrapped "A movie".
#   2. Appends 25 new paragraphs after it containing team meeting minutes
#      ("Minutes for 9/1/22", "Minutes for 9/6/22", "Minutes for 9/8/22",
#      plus blank spacer paragraphs).
#
# Word's object model does not expose a direct "remove this proofErr"
# call, so we rebuild the paragraph's XML without the markers, then use
# Range.InsertXML to splice in the new minutes paragraphs verbatim so the
# run layout (multiple <w:r> per paragraph in several places) matches the
# authored content exactly.

$d = $word.ActiveDocument

# --- Step 1: rebuild the last paragraph, dropping the <w:proofErr/> wrappers ---
$paraCount = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($paraCount)
$targetRange = $targetPara.Range

$movieParaXml = '<w:p><w:pPr><w:pStyle w:val="Normal"/></w:pPr><w:r><w:rPr/><w:t>A movie</w:t></w:r><w:r><w:rPr/><w:t xml:space="preserve"> themed application, where a user can type a movie, the app would then give feedback on where to stream, reviews, and other possibly informative </w:t></w:r><w:r><w:rPr/><w:t>variables</w:t></w:r><w:r><w:rPr/><w:t>.</w:t></w:r></w:p>'
$targetRange.InsertXML($movieParaXml)

# --- Step 2: create a fresh trailing paragraph, then replace its XML with ---
# --- all of the new "Minutes for ..." content (25 paragraphs total).     ---
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.InsertParagraphAfter() | Out-Null

$paraCount = $d.Paragraphs.Count
$insertTarget = $d.Paragraphs.Item($paraCount)

$newParasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t xml:space="preserve">Minutes for </w:t></w:r><w:r><w:rPr /><w:t>9/1/22</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Throughout</w:t></w:r><w:r><w:rPr /><w:t xml:space="preserve"> the class period:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>We decided the app idea we won’t to work on (movie app)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t xml:space="preserve">We talked about features our </w:t></w:r><w:r><w:rPr /><w:t>App could</w:t></w:r><w:r><w:rPr /><w:t xml:space="preserve"> have</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t xml:space="preserve">We </w:t></w:r><w:r><w:rPr /><w:t>decided on</w:t></w:r><w:r><w:rPr /><w:t xml:space="preserve"> the name of the app with the help of Dr. Kelley (CineSearch)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Bobby worked on the landing page for our app</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t xml:space="preserve">Spencer typed up our readme file for </w:t></w:r><w:r><w:rPr /><w:t>GitHub</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Isaiah worked on the logo and color scheme of app</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t xml:space="preserve">Omar watched tutorials on how to design web pages through </w:t></w:r><w:r><w:rPr /><w:t>GitHub</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Minutes for 9/6/22</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Most of class time was spent with the guest speaker</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>After the guest speaker finished, the team discussed designs for the box, as well as ideas for the pitch video</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Bobby showed the team how to ask him permission to push a repository on GitHub</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>We also decided that Bobby would monitor pushes and give permission to a push, as he will be spending the most time with the code and GitHub</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Minutes for 9/8/22</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Isaiah brought the box to class</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Isaiah and Spencer spent the time getting the proper materials to begin the box design</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Spencer began trimming our construction paper to cover box, while Isaiah got the logo ready while also setting up his laptop with GitHub and downloading files necessary for project</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr><w:r><w:rPr /><w:t>Bobby and Omar worked together further on the landing page</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Normal" /></w:pPr></w:p>'
$insertTarget.Range.InsertXML($newParasXml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
